$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value  = "to go (destination に/へ)"
$ws.Range("A3").Value  = "to go back; to return (destination に/へ)"
$ws.Range("A4").Value  = "to listen; to hear (～を)"
$ws.Range("A5").Value  = "to drink (～を)"
$ws.Range("A6").Value  = "to speak; to talk (language を/で)"
$ws.Range("A7").Value  = "to read (～を)"
$ws.Range("A52").Value = "to eat (～を)"
$ws.Range("A54").Value = "to see; to look at; to watch (～を)"
$ws.Range("A55").Value = "to come (destination に/へ)"
$ws.Range("A56").Value = "to do (～を)"
$ws.Range("A57").Value = "to study (～を)"
